$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.062.68"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.790.17"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'227.30"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'31.20"
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value = "'45.98"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("D10").Value = "'0.281"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").Value = "'0.0661"
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "2.046.39"
$ws.Range("D14").Value = "'11.46"
$ws.Range("E14").Value = "  +10.83%  "
$ws.Range("D15").Value = "1.787.28"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'0.636"
$ws.Range("E16").Value = "  -1.83%  "
$ws.Range("D17").Value = "34.054.50"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "'4.22"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").Value = "'69.62"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "'253.23"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").Value = "0.0₃0742"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'10.45"
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("D24").Value = "'4.30"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "'157.16"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("D27").Value = "'16.61"
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("D28").Value = "'7.03"
$ws.Range("E28").Value = "  -2.02%  "
$ws.Range("E29").Value = "  -2.66%  "
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").Value = "'3.83"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'0.0517"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").Value = "'3.62"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "1.452.93"
$ws.Range("E36").Value = "  -8.56%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "'0.632"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "'0.0187"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("D40").Value = "'83.49"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").Value = "'2.83"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'0.902"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("D45").Value = "'0.0511"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "1.945.96"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").Value = "'5.73"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D50").Value = "'11.81"
$ws.Range("E50").Value = "  +3.85%  "
$ws.Range("D51").Value = "'51.21"
$ws.Range("E51").Value = "  -4.88%  "
